$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Widen column B (target stored width is 43.42578125 "characters"; this
# runtime quantizes ColumnWidth to 1/6-character steps when it stores the
# <col> width, so feed it the value that lands closest to the target).
$ws.Columns.Item(2).ColumnWidth = 42.59244791666666

# Rows 37 and 38 don't exist yet in the used range (old dimension was A1:B36) -
# give them the same formatting (style + row height) as row 36 before filling
# them in, so the new rows don't end up with default formatting.
$ws.Range("A36:B36").Copy()
$ws.Range("A37:B38").PasteSpecial(-4122)
$ws.Rows.Item(37).RowHeight = 15.75
$ws.Rows.Item(38).RowHeight = 15.75

# New text values for rows 32-38 (column B)
$texts = @(
    "Im happy that I have done my homework quickly",
    "It was a sad day today because I lost my homework",
    "Ive lost my keys",
    "This project is a big deal to me",
    "Together its so much fun",
    "You cant do that",
    "Today was a good day"
)

# Row 32 carries its own (non-shared) copy of the audio3 formula, the same
# way row 17/row 2-3 started each of the earlier audio1/audio2 runs.
$ws.Range("A32").Formula = "=CONCATENATE(""/audio/audio3 ("", ROW() - 31, "").wav"")"

# Rows 33-38 are entered as one range assignment so they become a single
# shared-formula group (si="3"), mirroring rows A4:A11/A14:A16/A18:A31.
$ws.Range("A33:A38").Formula = "=CONCATENATE(""/audio/audio3 ("", ROW() - 31, "").wav"")"

for ($i = 0; $i -lt 7; $i++) {
    $row = 32 + $i
    $ws.Range("B$row").Value = $texts[$i]
}

# Scroll / selection state to match the final file
$ws.Range("B38").Select()
$excel.ActiveWindow.ScrollRow = 13
